$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new 2020 data point in column Q, matching the formatting of the
# existing 2019 column (P) for the header row and the data row.
$ws.Range("P4").Copy() | Out-Null
$ws.Range("Q4").PasteSpecial(-4122) | Out-Null
$ws.Range("Q4").Value = 2020

$ws.Range("P5").Copy() | Out-Null
$ws.Range("Q5").PasteSpecial(-4122) | Out-Null
$ws.Range("Q5").Value = 90.6

# Leave the selection where the author last left it when saving.
$ws.Range("P12").Select() | Out-Null
